# Apply the "flipkart demo test case" update to the workbook:
#  - add a new row to the "data" sheet with the Flipkart URL (reusing the
#    existing "visited link" style from A3)
#  - add a new "flipkart" worksheet (placed after "data") containing the
#    product name that the demo test will search for
#  - leave the selection/cursor the way it would end up after typing the
#    URL into A4 and hitting Enter (-> A5), and after typing the product
#    name into the new sheet's A1 and moving on to B1

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)

# --- "data" sheet: append the flipkart URL under the existing URLs ------
$data.Range("A4").Value = "https://www.flipkart.com/"

# Match the visited-link formatting already used on A3
$data.Range("A3").Copy()
$data.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- new "flipkart" sheet, inserted right after "data" ------------------
$flipkart = $wb.Worksheets.Add($null, $data)
$flipkart.Name = "flipkart"
$flipkart.Range("A1").Value = "ASUS Vivobook"
$flipkart.Range("B1").Select()

# Return focus to the data sheet, cursor parked below the new entry
$data.Range("A5").Select()
